# Update the "want to go" count (column F) on the "Exhibition" sheet (sheet 1)
# and the "All Types" sheet (sheet 4), matching the regenerated data snapshot.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 189
$ws1.Range("F6").Value = 3814
$ws1.Range("F8").Value = 115
$ws1.Range("F9").Value = 98
$ws1.Range("F10").Value = 82
$ws1.Range("F12").Value = 674
$ws1.Range("F14").Value = 938
$ws1.Range("F17").Value = 152
$ws1.Range("F18").Value = 65
$ws1.Range("F19").Value = 102
$ws1.Range("F21").Value = 3369
$ws1.Range("F22").Value = 5699
$ws1.Range("F24").Value = 22
$ws1.Range("F26").Value = 514
$ws1.Range("F28").Value = 3329
$ws1.Range("F29").Value = 347
$ws1.Range("F30").Value = 16
$ws1.Range("F31").Value = 2435
$ws1.Range("F35").Value = 194
$ws1.Range("F36").Value = 254
$ws1.Range("F37").Value = 344
$ws1.Range("F38").Value = 114
$ws1.Range("F43").Value = 33
$ws1.Range("F46").Value = 542

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 189
$ws4.Range("F6").Value = 3814
$ws4.Range("F8").Value = 115
$ws4.Range("F9").Value = 98
$ws4.Range("F11").Value = 82
$ws4.Range("F13").Value = 674
$ws4.Range("F15").Value = 938
$ws4.Range("F18").Value = 152
$ws4.Range("F19").Value = 65
$ws4.Range("F20").Value = 102
$ws4.Range("F22").Value = 3369
$ws4.Range("F23").Value = 5699
$ws4.Range("F25").Value = 22
$ws4.Range("F27").Value = 514
$ws4.Range("F29").Value = 3329
$ws4.Range("F30").Value = 347
$ws4.Range("F31").Value = 16
$ws4.Range("F32").Value = 2435
$ws4.Range("F36").Value = 194
$ws4.Range("F37").Value = 254
$ws4.Range("F38").Value = 344
$ws4.Range("F39").Value = 114
$ws4.Range("F44").Value = 33
$ws4.Range("F47").Value = 542

Write-Host "Applied F-column updates to sheets 1 and 4"
